$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'" + "28.430.26"
$ws.Cells.Item(2, 5).Value = "  +5.22%  "
$ws.Cells.Item(3, 4).Value = "'" + "1.815.19"
$ws.Cells.Item(3, 5).Value = "  +5.00%  "
$ws.Cells.Item(4, 4).Value = "'" + "0.9969"
$ws.Cells.Item(4, 5).Value = "  -0.36%  "
$ws.Cells.Item(5, 4).Value = "'" + "318.04"
$ws.Cells.Item(5, 5).Value = "  +2.43%  "
$ws.Cells.Item(6, 4).Value = "'" + "0.9970"
$ws.Cells.Item(6, 5).Value = "  -0.31%  "
$ws.Cells.Item(7, 4).Value = "'" + "0.5671"
$ws.Cells.Item(7, 5).Value = "  +16.63%  "
$ws.Cells.Item(8, 4).Value = "'" + "0.3839"
$ws.Cells.Item(8, 5).Value = "  +10.09%  "
$ws.Cells.Item(9, 4).Value = "'" + "43.44"
$ws.Cells.Item(9, 5).Value = "  +0.15%  "
$ws.Cells.Item(10, 4).Value = "'" + "0.07635"
$ws.Cells.Item(10, 5).Value = "  +5.42%  "
$ws.Cells.Item(11, 4).Value = "'" + "1.139"
$ws.Cells.Item(11, 5).Value = "  +8.11%  "
$ws.Cells.Item(12, 4).Value = "'" + "21.37"
$ws.Cells.Item(12, 5).Value = "  +6.96%  "
$ws.Cells.Item(13, 4).Value = "'" + "0.9968"
$ws.Cells.Item(13, 5).Value = "  -0.38%  "
$ws.Cells.Item(14, 4).Value = "'" + "6.241"
$ws.Cells.Item(15, 4).Value = "'" + "1.804.36"
$ws.Cells.Item(15, 5).Value = "  +4.40%  "
$ws.Cells.Item(16, 4).Value = "'" + "7.255"
$ws.Cells.Item(16, 5).Value = "  +5.78%  "
$ws.Cells.Item(17, 4).Value = "'" + "92.33"
$ws.Cells.Item(17, 5).Value = "  +6.11%  "
$ws.Cells.Item(18, 4).Value = "'" + "0.00001081"
$ws.Cells.Item(18, 5).Value = "  +4.56%  "
$ws.Cells.Item(19, 4).Value = "'" + "0.06520"
$ws.Cells.Item(19, 5).Value = "  +1.96%  "
$ws.Cells.Item(20, 4).Value = "'" + "0.9970"
$ws.Cells.Item(20, 5).Value = "  -0.33%  "
$ws.Cells.Item(21, 4).Value = "'" + "17.30"
$ws.Cells.Item(21, 5).Value = "  +3.97%  "
$ws.Cells.Item(22, 4).Value = "'" + "6.005"
$ws.Cells.Item(22, 5).Value = "  +4.92%  "
$ws.Cells.Item(23, 4).Value = "'" + "28.428.70"
$ws.Cells.Item(23, 5).Value = "  +4.97%  "
$ws.Cells.Item(24, 5).Value = "  +3.09%  "
$ws.Cells.Item(25, 4).Value = "'" + "2.103"
$ws.Cells.Item(25, 5).Value = "  +1.42%  "
$ws.Cells.Item(26, 4).Value = "'" + "20.86"
$ws.Cells.Item(26, 5).Value = "  +4.20%  "
$ws.Cells.Item(27, 4).Value = "'" + "156.90"
$ws.Cells.Item(28, 4).Value = "'" + "2.391"
$ws.Cells.Item(28, 5).Value = "  +15.25%  "
$ws.Cells.Item(29, 4).Value = "'" + "2.015.88"
$ws.Cells.Item(29, 5).Value = "  +4.70%  "
$ws.Cells.Item(30, 4).Value = "'" + "123.70"
$ws.Cells.Item(30, 5).Value = "  +2.23%  "
$ws.Cells.Item(31, 4).Value = "'" + "1.148"
$ws.Cells.Item(31, 5).Value = "  +9.86%  "
$ws.Cells.Item(32, 4).Value = "'" + "0.1051"
$ws.Cells.Item(32, 5).Value = "  +12.42%  "
$ws.Cells.Item(33, 4).Value = "'" + "5.777"
$ws.Cells.Item(33, 5).Value = "  +7.00%  "
$ws.Cells.Item(34, 4).Value = "'" + "3.626"
$ws.Cells.Item(34, 5).Value = "  -0.51%  "
$ws.Cells.Item(35, 4).Value = "'" + "0.02319"
$ws.Cells.Item(35, 5).Value = "  +6.14%  "
$ws.Cells.Item(36, 4).Value = "'" + "0.2140"
$ws.Cells.Item(36, 5).Value = "  +7.18%  "
$ws.Cells.Item(37, 4).Value = "'" + "8.724"
$ws.Cells.Item(37, 5).Value = "  +15.81%  "
$ws.Cells.Item(38, 4).Value = "'" + "11.68"
$ws.Cells.Item(38, 5).Value = "  +6.28%  "
$ws.Cells.Item(39, 2).Value = "TheSandbox"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(39, 4).Value = "'" + "0.6438"
$ws.Cells.Item(39, 5).Value = "  +7.47%  "
$ws.Cells.Item(40, 2).Value = "Hedera"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(40, 4).Value = "'" + "0.06089"
$ws.Cells.Item(40, 5).Value = "  +2.73%  "
$ws.Cells.Item(41, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(41, 4).Value = "'" + "5.048"
$ws.Cells.Item(41, 5).Value = "  +6.07%  "
$ws.Cells.Item(42, 4).Value = "'" + "0.9968"
$ws.Cells.Item(42, 5).Value = "  -0.29%  "
$ws.Cells.Item(43, 4).Value = "'" + "1.156"
$ws.Cells.Item(43, 5).Value = "  +3.36%  "
$ws.Cells.Item(44, 4).Value = "'" + "1.376"
$ws.Cells.Item(44, 5).Value = "  -3.43%  "
$ws.Cells.Item(45, 5).Value = "  +5.06%  "
$ws.Cells.Item(46, 4).Value = "'" + "0.6010"
$ws.Cells.Item(46, 5).Value = "  +6.77%  "
$ws.Cells.Item(47, 4).Value = "'" + "3.699"
$ws.Cells.Item(47, 5).Value = "  +3.31%  "
$ws.Cells.Item(48, 4).Value = "'" + "122.56"
$ws.Cells.Item(48, 5).Value = "  +2.95%  "
$ws.Cells.Item(49, 4).Value = "'" + "1.942"
$ws.Cells.Item(49, 5).Value = "  +4.92%  "
$ws.Cells.Item(50, 4).Value = "'" + "1.146"
$ws.Cells.Item(50, 5).Value = "  +3.61%  "
$ws.Cells.Item(51, 4).Value = "'" + "0.06842"
$ws.Cells.Item(51, 5).Value = "  +2.97%  "
